# Apply finalized SeenRx CKJ report values (rows 2-36, columns B-J)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2=37, C2=19, D2=20, E2=4, F2=5, G2=25, H2=36, I2=10, J2=42
$ws.Range("B2").Value = 37
$ws.Range("C2").Value = 19
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = 36
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 42

# Row 3: B3=7, C3=7, D3=8, E3=0, F3=1, G3=7, H3=6, I3=2, J3=7
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 7

# Row 4: B4=1, C4=1, D4=3, E4=0, F4=0, G4=1, H4=1, I4=1, J4=1
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1

# Row 5: B5=4, C5=5, D5=1, E5=0, F5=0, G5=1, H5=1, I5=0, J5=2
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2

# Row 6: B6=0, C6=0, D6=3, E6=0, F6=0, G6=2, H6=1, I6=1, J6=2
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2

# Row 7: B7=2, C7=0, D7=1, E7=0, F7=0, G7=0, H7=0, I7=0, J7=0
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0

# Row 8: B8=0, C8=1, D8=0, E8=0, F8=1, G8=3, H8=3, I8=0, J8=2
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2

# Row 9: B9=0, C9=0, D9=0, E9=0, F9=0, G9=0, H9=0, I9=0, J9=0
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0

# Row 10: B10=11, C10=7, D10=2, E10=2, F10=3, G10=6, H10=14, I10=1, J10=3
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 6
$ws.Range("H10").Value = 14
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 3

# Row 11: B11=0, C11=0, D11=0, E11=0, F11=0, G11=0, H11=0, I11=0, J11=0
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0

# Row 12: B12=5, C12=4, D12=0, E12=2, F12=0, G12=0, H12=2, I12=0, J12=0
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0

# Row 13: B13=0, C13=0, D13=0, E13=0, F13=0, G13=3, H13=5, I13=0, J13=1
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1

# Row 14: B14=6, C14=2, D14=1, E14=0, F14=3, G14=0, H14=5, I14=0, J14=1
$ws.Range("B14").Value = 6
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1

# Row 15: B15=0, C15=0, D15=0, E15=0, F15=0, G15=3, H15=0, I15=1, J15=1
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1

# Row 16: B16=0, C16=1, D16=1, E16=0, F16=0, G16=0, H16=2, I16=0, J16=0
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0

# Row 17: B17=5, C17=3, D17=6, E17=0, F17=0, G17=3, H17=8, I17=4, J17=16
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 8
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 16

# Row 18: B18=1, C18=3, D18=0, E18=0, F18=0, G18=0, H18=1, I18=0, J18=1
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 1

# Row 19: B19=0, C19=0, D19=2, E19=0, F19=0, G19=0, H19=5, I19=2, J19=1
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 2
$ws.Range("J19").Value = 1

# Row 20: B20=4, C20=0, D20=0, E20=0, F20=0, G20=2, H20=1, I20=0, J20=0
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0

# Row 21: B21=0, C21=0, D21=0, E21=0, F21=0, G21=0, H21=0, I21=2, J21=8
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = 8

# Row 22: B22=0, C22=0, D22=4, E22=0, F22=0, G22=1, H22=1, I22=0, J22=6
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 6

# Row 23: B23=5, C23=0, D23=2, E23=0, F23=0, G23=4, H23=3, I23=2, J23=14
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = 14

# Row 24: B24=0, C24=0, D24=0, E24=0, F24=0, G24=0, H24=0, I24=0, J24=0
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0

# Row 25: B25=2, C25=0, D25=0, E25=0, F25=0, G25=3, H25=3, I25=2, J25=13
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 13

# Row 26: B26=0, C26=0, D26=0, E26=0, F26=0, G26=0, H26=0, I26=0, J26=0
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0

# Row 27: B27=3, C27=0, D27=0, E27=0, F27=0, G27=0, H27=0, I27=0, J27=0
$ws.Range("B27").Value = 3
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0

# Row 28: B28=0, C28=0, D28=0, E28=0, F28=0, G28=0, H28=0, I28=0, J28=0
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0

# Row 29: B29=0, C29=0, D29=2, E29=0, F29=0, G29=1, H29=0, I29=0, J29=1
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1

# Row 30: B30=9, C30=2, D30=2, E30=2, F30=1, G30=5, H30=5, I30=1, J30=2
$ws.Range("B30").Value = 9
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 2
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 5
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 2

# Row 31: B31=0, C31=0, D31=2, E31=0, F31=0, G31=2, H31=1, I31=0, J31=1
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 1
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1

# Row 32: B32=2, C32=0, D32=0, E32=0, F32=0, G32=0, H32=0, I32=0, J32=0
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0

# Row 33: B33=2, C33=0, D33=0, E33=0, F33=1, G33=2, H33=0, I33=1, J33=1
$ws.Range("B33").Value = 2
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = 1

# Row 34: B34=0, C34=1, D34=0, E34=2, F34=0, G34=0, H34=1, I34=0, J34=0
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 2
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 1
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0

# Row 35: B35=1, C35=1, D35=0, E35=0, F35=0, G35=0, H35=1, I35=0, J35=0
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 1
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0

# Row 36: B36=4, C36=0, D36=0, E36=0, F36=0, G36=2, H36=2, I36=0, J36=0
$ws.Range("B36").Value = 4
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 2
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0

